# Applies the "Updated cryptos list" data refresh to the crypto table.
# Columns: A=rank(unchanged) B=Coin C=Link D=Price E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.775.09"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "3.413.73"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'578.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").Value = "'143.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.473"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'7.65"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").Value = "'0.123"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").Value = "'0.386"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "3.985.07"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").Value = "'28.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "3.403.26"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "'0.0000170"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").Value = "61.788.12"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "'6.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "'13.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").Value = "'9.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.80%  "
$ws.Range("D21").Value = "'388.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").Value = "'74.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").Value = "'0.551"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'0.0000115"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "'7.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("D29").Value = "'8.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("D30").Value = "'2.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'23.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D34").Value = "'6.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").Value = "'5.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.80%  "
$ws.Range("D36").Value = "'168.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("D37").Value = "3.442.78"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").Value = "'1.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "'28.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.59%  "
$ws.Range("D40").Value = "'0.0759"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("D41").Value = "'0.785"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").Value = "'4.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("D43").Value = "'1.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("D44").Value = "'1.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.62%  "
$ws.Range("D45").Value = "2.500.26"
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("D46").Value = "'22.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").Value = "'6.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "'0.0264"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'2.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.15%  "
$ws.Range("D51").Value = "'0.207"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.28%  "
